# edit.ps1 - reproduces the tracked content changes from the diff:
#   1. Remove the old "_GoBack" bookmark that wrapped "136-Я".
#   2. Merge the spell-checked "Студента (" / "ки" / ")" runs into one run
#      "Студента (ки)" (drops the now-stale proofErr spell-check markers).
#   3. Expand "МГГТК АГУ" -> "МГГТК ФГБОУ ВО АГУ", split across three runs
#      (matching formatting), and park a fresh collapsed "_GoBack" bookmark
#      between the 2nd and 3rd run (where the edit ended).
#   4. Merge the grammar-checked "Проведение анализа рисков и характеристик "
#      / "качества " / " ПО" / " " runs into two runs (drops the proofErr
#      grammar-check markers), keeping the exact same text.

$d = $word.ActiveDocument

# --- 1. Drop the stale _GoBack bookmark around "136-Я" ---------------------
$d.Bookmarks("_GoBack").Delete()

# --- 2. Merge "Студента (" + "ки" + ")" into a single run ------------------
$studentTarget = "Студента (ки)"
$d.Content.Find.Execute($studentTarget, $true, $false, $false, $false, $false, `
    $true, 1, $false, $studentTarget, 2) | Out-Null

# --- 3. "МГГТК АГУ" -> "МГГТК" / " ФГБОУ ВО" / " АГУ" (+ new _GoBack) ------
$d.Content.Find.Execute("МГГТК АГУ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "МГГТК ФГБОУ ВО АГУ", 2) | Out-Null

# Locate the boundary right after "МГГТК" (start of the inserted text).
$afterMggtk = $d.Content
$afterMggtk.Find.Execute("МГГТК", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$splitPoint1 = $afterMggtk.End

# Locate the boundary right after "ФГБОУ ВО" (end of the inserted text /
# where the user's cursor would have landed -> becomes the new _GoBack spot).
$afterFgbou = $d.Content
$afterFgbou.Find.Execute("ФГБОУ ВО", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$splitPoint2 = $afterFgbou.End

# Locate the boundary right after " АГУ" (before the pre-existing, untouched
# run of trailing spaces that follows it in the same cell).
$afterAgu = $d.Content
$afterAgu.Find.Execute(" АГУ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$splitPoint3 = $afterAgu.End

# Toggling a character property off/on over a sub-range forces the engine to
# split the run exactly at those boundaries without altering the visible
# formatting (since we flip it right back to its original value).
$middle = $d.Range($splitPoint1, $splitPoint2)
$middle.Font.Bold = 1
$middle.Font.Bold = 0

$tail = $d.Range($splitPoint2, $splitPoint3)
$tail.Font.Bold = 1
$tail.Font.Bold = 0

# Drop a fresh, collapsed _GoBack bookmark at the end of the inserted text;
# adding it also splits the run there, giving the 3rd run (" АГУ").
$goBackRange = $d.Range($splitPoint2, $splitPoint2)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

# --- 4. Merge the "Проведение анализа рисков..." runs into two runs --------
$riskTarget = "Проведение анализа рисков и характеристик качества  ПО "
$riskFull = $d.Content
$riskFull.Find.Execute($riskTarget, $true, $false, $false, $false, $false, `
    $true, 1, $false, $riskTarget, 2) | Out-Null
$riskFullEnd = $riskFull.End

$riskHead = $d.Content
$riskHead.Find.Execute("Проведение анализа рисков и характеристик качества ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$riskSplit = $riskHead.End

$riskTail = $d.Range($riskSplit, $riskFullEnd)
$riskTail.Font.Bold = 1
$riskTail.Font.Bold = 0
